# bondxll/bondlib.xlsx — "date::dpy conforms with chrono move date enumes to xll_date.cpp"
#
# The real commit re-wired the per-row bond formulas to reference the row's own
# B/C/D cells directly instead of the whole-column named ranges (maturity/coupon/bond),
# dropped the scratch date-math cells in row 4 (G4/H4) and the per-row TMX.INSTRUMENT
# error column (F7), and added one spilling TMX.INSTRUMENT call anchored at F16 that
# fans out across F16:Y17. None of the _xll.* functions are resolvable in this engine
# (no XLL host registered), so their cached results cannot be recomputed here — we
# reproduce the formula text/shape exactly and leave those specific cells' cached
# values as Excel would leave an unresolved UDF, while every literal (non-_xll) value
# is written verbatim.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 4: drop the scratch "dated + 1 year" / "year-fraction" helper cells.
# G4 keeps its date-formatted style but loses its formula+value; H4 disappears
# entirely (ClearContents on a cell with no explicit style drops the <c> node).
# ---------------------------------------------------------------------------
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()

# ---------------------------------------------------------------------------
# B2 ("BondLib" banner) picks up center/middle alignment on top of its existing
# large font.
# ---------------------------------------------------------------------------
$ws.Range("B2").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B2").VerticalAlignment = -4108     # xlCenter

# ---------------------------------------------------------------------------
# Rows 7-16: re-point the bond/instrument formulas at the row's own cells
# instead of the maturity/coupon/bond named ranges, and give column E the same
# look as column D (the custom "0x#" number format + small red font + centered).
# FormulaArray (rather than Formula) keeps these as legacy array formulas with
# an explicit single-cell ref="..." like the original file, and — since the
# _xll functions can't be evaluated here — leaves the cell without a bogus
# cached value instead of baking in a spurious #NAME?.
# ---------------------------------------------------------------------------
for ($r = 7; $r -le 16; $r++) {
    $ws.Range("D$r").FormulaArray = "=_xll.\TMX.BOND.SIMPLE(B$r, C$r)"

    $ws.Range("E$r").NumberFormat = $ws.Range("D$r").NumberFormat
    $ws.Range("E$r").Font.Size = $ws.Range("D$r").Font.Size
    $ws.Range("E$r").Font.Color = $ws.Range("D$r").Font.Color
    $ws.Range("E$r").HorizontalAlignment = $ws.Range("D$r").HorizontalAlignment
    $ws.Range("E$r").VerticalAlignment = $ws.Range("D$r").VerticalAlignment
    $ws.Range("E$r").FormulaArray = "=_xll.\TMX.BOND.INSTRUMENT(D$r, dated)"
}

# Row 7's old per-row TMX.INSTRUMENT error cell is gone.
$ws.Range("F7").ClearContents()

# ---------------------------------------------------------------------------
# Row 16/17: the new spilling TMX.INSTRUMENT(E16) call anchored at F16, fanning
# out across F16:Y17. Seed the whole block with the literal spill results first
# (Excel won't let you touch part of an array afterwards), then drop the array
# formula onto the same range so the anchor cell carries
# f t="array" ref="F16:Y17" while the rest of the block keeps its values —
# exactly how the saved workbook represents a spilled legacy array formula.
# ---------------------------------------------------------------------------
$row16 = @(0.9965981505438167, 1.4948972258157252, 1.9931963010876335, 2.4914953763595418, 2.9897944516314503, 3.488093526903359, 3.986392602175267, 4.4846916774471755, 4.9829907527190835, 5.4812898279909925, 5.979588903262901, 6.477887978534809, 6.976187053806718, 7.474486129078626, 7.972785204350534, 8.471084279622442, 8.969383354894351, 9.46768243016626, 9.965981505438167)
$row17 = @(0.02486111111111111, 0.02486111111111111, 0.02486111111111111, 0.025, 0.02486111111111111, 0.025, 0.02486111111111111, 0.02513888888888889, 0.024722222222222225, 0.025, 0.024722222222222225, 0.02513888888888889, 0.024722222222222225, 0.02513888888888889, 0.024722222222222225, 0.02513888888888889, 0.024722222222222225, 0.025, 0.024722222222222225, 0.02513888888888889)

$spill = New-Object 'object[,]' 2,20
for ($c = 0; $c -lt 20; $c++) {
    $spill[0, $c] = $row16[$c]
    $spill[1, $c] = $row17[$c]
}
$ws.Range("F16:Y17").Value = $spill
$ws.Range("F16:Y17").FormulaArray = "=_xll.TMX.INSTRUMENT(E16)"

# ---------------------------------------------------------------------------
# Selection follows the edited area.
# ---------------------------------------------------------------------------
$ws.Range("E7").Select()
